$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link / Volume(1h) columns (plain text assignment)
$ws.Range('E2').Value = '  -0.09%  '
$ws.Range('E3').Value = '  +0.05%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('E5').Value = '  +2.07%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('E7').Value = '  -0.85%  '
$ws.Range('E8').Value = '  -1.94%  '
$ws.Range('B9').Value = 'Dogecoin'
$ws.Range('C9').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('B10').Value = 'WrappedEther'
$ws.Range('C10').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('E10').Value = '  +0.40%  '
$ws.Range('B11').Value = 'TRON'
$ws.Range('C11').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E11').Value = '  +1.88%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('E12').Value = '  -1.86%  '
$ws.Range('B13').Value = 'Polygon'
$ws.Range('C13').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('E13').Value = '  -1.33%  '
$ws.Range('B14').Value = 'Polkadot'
$ws.Range('C14').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('B15').Value = 'Litecoin'
$ws.Range('C15').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E15').Value = '  +0.26%  '
$ws.Range('B16').Value = 'Dai'
$ws.Range('C16').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('E16').Value = '  +0.02%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('E17').Value = '  +0.10%  '
$ws.Range('B18').Value = 'WrappedBTC'
$ws.Range('C18').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('E18').Value = '  -0.04%  '
$ws.Range('B19').Value = 'Avalanche'
$ws.Range('C19').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('E19').Value = '  +1.97%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('E20').Value = '  +1.08%  '
$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('E22').Value = '  +2.52%  '
$ws.Range('B23').Value = 'Cosmos'
$ws.Range('C23').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('E23').Value = '  -1.26%  '
$ws.Range('B24').Value = 'Chainlink'
$ws.Range('C24').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('E24').Value = '  +1.73%  '
$ws.Range('B25').Value = 'Monero'
$ws.Range('C25').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E25').Value = '  -0.45%  '
$ws.Range('B26').Value = 'Toncoin'
$ws.Range('C26').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('B27').Value = 'EthereumClassic'
$ws.Range('C27').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('E27').Value = '  +0.49%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('E28').Value = '  -1.44%  '
$ws.Range('B29').Value = 'BitcoinCash'
$ws.Range('C29').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('E29').Value = '  +3.09%  '
$ws.Range('B30').Value = 'InternetComputer(DFINITY)'
$ws.Range('C30').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('E30').Value = '  +5.83%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('B32').Value = 'Filecoin'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('E32').Value = '  +5.16%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('B34').Value = 'HuobiToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('E34').Value = '  +0.91%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('E35').Value = '  +0.41%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('B37').Value = 'MXToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('E37').Value = '  +1.45%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E38').Value = '  +0.61%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('E39').Value = '  -1.00%  '
$ws.Range('E40').Value = '  +0.07%  '
$ws.Range('B41').Value = 'Quant'
$ws.Range('C41').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('E41').Value = '  -3.06%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('E42').Value = '  +0.49%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('E43').Value = '  +1.27%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('E44').Value = '  -0.72%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('E45').Value = '  +0.10%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('E46').Value = '  -2.06%  '
$ws.Range('B47').Value = 'Aptos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('E47').Value = '  -1.14%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('E48').Value = '  +2.90%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('E49').Value = '  +0.97%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E50').Value = '  -0.51%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('E51').Value = '  +0.90%  '

# Update Price column values.
# Force the whole Price column to Text format first so strings such as
# "1.000", "15.10" or "0.000006690" are stored verbatim instead of being
# re-interpreted (and visually truncated) as numbers by Excel.
$ws.Range('D2:D51').NumberFormat = '@'
$ws.Range('D2').Value = '25.819.48'
$ws.Range('D3').Value = '1.734.67'
$ws.Range('D5').Value = '236.11'
$ws.Range('D6').Value = '1.001'
$ws.Range('D7').Value = '0.5116'
$ws.Range('D8').Value = '0.2719'
$ws.Range('D9').Value = '0.06101'
$ws.Range('D10').Value = '1.742.11'
$ws.Range('D11').Value = '0.07164'
$ws.Range('D12').Value = '14.91'
$ws.Range('D13').Value = '0.6319'
$ws.Range('D14').Value = '4.578'
$ws.Range('D15').Value = '76.96'
$ws.Range('D16').Value = '1.001'
$ws.Range('D17').Value = '1.001'
$ws.Range('D18').Value = '25.823.65'
$ws.Range('D19').Value = '11.69'
$ws.Range('D20').Value = '0.000006690'
$ws.Range('D21').Value = '1.965.51'
$ws.Range('D22').Value = '4.237'
$ws.Range('D23').Value = '8.607'
$ws.Range('D24').Value = '5.207'
$ws.Range('D25').Value = '138.97'
$ws.Range('D26').Value = '1.508'
$ws.Range('D27').Value = '15.10'
$ws.Range('D28').Value = '1.750'
$ws.Range('D29').Value = '105.11'
$ws.Range('D30').Value = '3.896'
$ws.Range('D31').Value = '0.08348'
$ws.Range('D32').Value = '3.600'
$ws.Range('D33').Value = '0.04544'
$ws.Range('D34').Value = '2.643'
$ws.Range('D35').Value = '0.9787'
$ws.Range('D36').Value = '0.6187'
$ws.Range('D37').Value = '2.687'
$ws.Range('D38').Value = '0.01587'
$ws.Range('D39').Value = '1.911'
$ws.Range('D40').Value = '1.000'
$ws.Range('D41').Value = '97.38'
$ws.Range('D42').Value = '0.3827'
$ws.Range('D43').Value = '0.7321'
$ws.Range('D44').Value = '4.925'
$ws.Range('D45').Value = '0.1123'
$ws.Range('D46').Value = '0.05265'
$ws.Range('D47').Value = '6.156'
$ws.Range('D48').Value = '54.43'
$ws.Range('D49').Value = '30.32'
$ws.Range('D50').Value = '7.548'
$ws.Range('D51').Value = '0.3397'

# Restore the original (default) cell formatting now that the text is locked in,
# so no stray number-format style is left behind on these cells.
$ws.Range('D2:D51').ClearFormats()
